$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'Purbaya Janji Bayar Tunggakan Rp 55 T ke BUMN: Tapi Jangan Rugi Terus!'
$ws.Range("B2").Value = '2025-09-30T14:03:17+07:00'
$ws.Range("C2").Value = 'Anisa Indraini'
$ws.Range("D2").Value = 'https://finance.detik.com/berita-ekonomi-bisnis/d-8137535/purbaya-janji-bayar-tunggakan-rp-55-t-ke-bumn-tapi-jangan-rugi-terus'

$ws.Range("A3").Value = 'Purbaya Bantah Dikte Bank BUMN Naikkan Bunga Deposito Valas ke 4%'
$ws.Range("B3").Value = '2025-09-26T16:11:28+07:00'
$ws.Range("C3").Value = 'Anisa Indraini'
$ws.Range("D3").Value = 'https://finance.detik.com/moneter/d-8131959/purbaya-bantah-dikte-bank-bumn-naikkan-bunga-deposito-valas-ke-4'

$ws.Range("A4").Value = 'Purbaya Beri Subsidi Bunga buat KUR Perumahan, Segini Besarannya'
$ws.Range("B4").Value = '2025-09-26T11:00:52+07:00'
$ws.Range("C4").Value = 'Danica Adhitiawarman'
$ws.Range("D4").Value = 'https://www.detik.com/properti/berita/d-8131197/purbaya-beri-subsidi-bunga-buat-kur-perumahan-segini-besarannya'

$ws.Range("A5").Value = 'Janji Purbaya Genjot Pertumbuhan Ekonomi-Tekan Jumlah Utang'
$ws.Range("B5").Value = '2025-09-24T08:58:32+07:00'
$ws.Range("C5").Value = 'Anisa Indraini'
$ws.Range("D5").Value = 'https://finance.detik.com/berita-ekonomi-bisnis/d-8126940/janji-purbaya-genjot-pertumbuhan-ekonomi-tekan-jumlah-utang'

$ws.Range("A6").Value = 'Purbaya: Saya Tak Akan Tambah Utang Terlalu Besar'
$ws.Range("B6").Value = '2025-09-23T13:13:19+07:00'
$ws.Range("C6").Value = 'Anisa Indraini'
$ws.Range("D6").Value = 'https://finance.detik.com/berita-ekonomi-bisnis/d-8125951/purbaya-saya-tak-akan-tambah-utang-terlalu-besar'

$ws.Range("A7").Value = 'Purbaya Ingatkan Bank Penerima Rp 200 T: Kalau Nggak Hati-hati Bisa Dipecat'
$ws.Range("B7").Value = '2025-09-16T18:29:51+07:00'
$ws.Range("C7").Value = 'Herdi Alif Al Hikam'
$ws.Range("D7").Value = 'https://finance.detik.com/moneter/d-8115243/purbaya-ingatkan-bank-penerima-rp-200-t-kalau-nggak-hati-hati-bisa-dipecat'

$ws.Range("A8").Value = 'Purbaya Sidak Kantor BNI Saat Direksi Sedang Rapat, Ada Apa?'
$ws.Range("B8").Value = '2025-09-29T13:04:09+07:00'
$ws.Range("C8").Value = 'Anisa Indraini'
$ws.Range("D8").Value = 'https://finance.detik.com/moneter/d-8135534/purbaya-sidak-kantor-bni-saat-direksi-sedang-rapat-ada-apa'

$ws.Range("A9").Value = 'Ulama Jatim Dukung Menkeu Berantas Rokok Ilegal'
$ws.Range("B9").Value = '2025-09-30T15:10:55+07:00'
$ws.Range("C9").Value = 'Faiq Azmi'
$ws.Range("D9").Value = 'https://www.detik.com/jatim/berita/d-8137673/ulama-jatim-dukung-menkeu-berantas-rokok-ilegal'

$ws.Range("A10").Value = 'Usai Bertemu Djarum-Gudang Garam, Purbaya Putuskan Cukai Rokok Tak Naik'
$ws.Range("B10").Value = '2025-09-27T06:27:00+07:00'
$ws.Range("C10").Value = 'Anisa Indraini'
$ws.Range("D10").Value = 'https://finance.detik.com/industri/d-8132587/usai-bertemu-djarum-gudang-garam-purbaya-putuskan-cukai-rokok-tak-naik'

$ws.Range("A11").Value = 'Kapal Baharkam Polri Gagalkan Penyelundupan Rokok Ilegal dari Batam ke Riau'
$ws.Range("B11").Value = '2025-09-26T20:21:28+07:00'
$ws.Range("C11").Value = 'Alamudin Hamapu'
$ws.Range("D11").Value = 'https://www.detik.com/sumut/hukum-dan-kriminal/d-8132436/kapal-baharkam-polri-gagalkan-penyelundupan-rokok-ilegal-dari-batam-ke-riau'

$ws.Range("A12").Value = 'Terungkap Modus Jualan Rokok Ilegal di Toko Online'
$ws.Range("B12").Value = '2025-09-26T19:17:54+07:00'
$ws.Range("C12").Value = 'Amanda Christabel'
$ws.Range("D12").Value = 'https://finance.detik.com/industri/d-8132372/terungkap-modus-jualan-rokok-ilegal-di-toko-online'

$ws.Range("A13").Value = 'Jabar Dikepung Rokok Ilegal, Bikin Negara Boncos Miliaran'
$ws.Range("B13").Value = '2025-09-26T14:22:03+07:00'
$ws.Range("C13").Value = 'Rifat Alhamidi'
$ws.Range("D13").Value = 'https://www.detik.com/jabar/berita/d-8131643/jabar-dikepung-rokok-ilegal-bikin-negara-boncos-miliaran'

$ws.Range("A14").Value = 'Bea Cukai Kupang Sita 58.165 Batang Rokok Ilegal'
$ws.Range("B14").Value = '2025-09-25T21:52:39+07:00'
$ws.Range("C14").Value = 'Simon Selly'
$ws.Range("D14").Value = 'https://www.detik.com/bali/hukum-dan-kriminal/d-8130829/bea-cukai-kupang-sita-58-165-batang-rokok-ilegal'

$ws.Range("A15").Value = 'Bea Cukai Musnahkan 64 Juta Rokok dan MMEA Ilegal'
$ws.Range("B15").Value = '2025-06-25T22:30:31+07:00'
$ws.Range("C15").Value = '-'
$ws.Range("D15").Value = 'https://finance.detik.com/foto-bisnis/d-7982307/bea-cukai-musnahkan-64-juta-rokok-dan-mmea-ilegal'

$ws.Range("A16").Value = '20 Juta Batang Rokok Ilegal Bea Cukai Dimusnahkan'
$ws.Range("B16").Value = '2024-11-21T20:35:11+07:00'
$ws.Range("C16").Value = '-'
$ws.Range("D16").Value = 'https://news.detik.com/foto-news/d-7649790/20-juta-batang-rokok-ilegal-bea-cukai-dimusnahkan'

$ws.Range("A17").Value = 'Pengungkapan Jutaan Rokok Ilegal di Surabaya'
$ws.Range("B17").Value = '2024-11-11T19:33:27+07:00'
$ws.Range("C17").Value = '-'
$ws.Range("D17").Value = 'https://www.detik.com/jatim/foto/d-7633599/pengungkapan-jutaan-rokok-ilegal-di-surabaya'

$ws.Range("A18").Value = 'Potret Bea Cukai Bekasi Musnahkan 5 Juta Batang Rokok Ilegal'
$ws.Range("B18").Value = '2024-10-09T17:40:50+07:00'
$ws.Range("C18").Value = '-'
$ws.Range("D18").Value = 'https://finance.detik.com/foto-bisnis/d-7580101/potret-bea-cukai-bekasi-musnahkan-5-juta-batang-rokok-ilegal'

$ws.Range("A19").Value = 'Bea Cukai Musnahkan Rokok hingga Miras Ilegal di Semarang'
$ws.Range("B19").Value = '2024-07-09T19:00:38+07:00'
$ws.Range("C19").Value = '-'
$ws.Range("D19").Value = 'https://news.detik.com/foto-news/d-7430343/bea-cukai-musnahkan-rokok-hingga-miras-ilegal-di-semarang'

$ws.Range("A20").Value = 'Anggota DPR Usulkan Moratorium Cukai Rokok 3 Tahun'
$ws.Range("B20").Value = '2025-09-25T11:26:57+07:00'
$ws.Range("C20").Value = 'Andi Hidayat'
$ws.Range("D20").Value = 'https://finance.detik.com/industri/d-8129626/anggota-dpr-usulkan-moratorium-cukai-rokok-3-tahun'

$ws.Range("A21").Value = 'Pengedar Rokok Ilegal Kalsel Ngaku Untung Rp 500 Per Bungkus'
$ws.Range("B21").Value = '2025-09-23T20:30:16+07:00'
$ws.Range("C21").Value = 'Khairun Nisa'
$ws.Range("D21").Value = 'https://www.detik.com/kalimantan/hukum-dan-kriminal/d-8126582/pengedar-rokok-ilegal-kalsel-ngaku-untung-rp-500-per-bungkus'

$ws.Range("A22").Value = 'Ditpolairud Polda Kalsel Ungkap Kasus Cukai Rokok Rugikan Negara Rp 505 Juta'
$ws.Range("B22").Value = '2025-09-23T16:45:36+07:00'
$ws.Range("C22").Value = 'Khairun Nisa'
$ws.Range("D22").Value = 'https://www.detik.com/kalimantan/hukum-dan-kriminal/d-8126392/ditpolairud-polda-kalsel-ungkap-kasus-cukai-rokok-rugikan-negara-rp-505-juta'

$ws.Range("A23").Value = 'Cukai Rokok Tinggi Dinilai Picu Rokok Ilegal Merajalela'
$ws.Range("B23").Value = '2025-09-23T13:02:25+07:00'
$ws.Range("C23").Value = 'Andi Hidayat'
$ws.Range("D23").Value = 'https://finance.detik.com/industri/d-8126097/cukai-rokok-tinggi-dinilai-picu-rokok-ilegal-merajalela'
